$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 493; this shifts existing rows 493:563 down to 494:564
$ws.Range("A493:R493").EntireRow.Insert()

# Populate the newly inserted row 493 with the new record's data
$ws.Cells.Item(493, 1).Value = 5
$ws.Cells.Item(493, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(493, 3).Value = "Maule"
$ws.Cells.Item(493, 4).Value = 45077
$ws.Cells.Item(493, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(493, 5).Value = 7
$ws.Cells.Item(493, 6).Value = 100112032
$ws.Cells.Item(493, 7).Value = "Zapallo italiano"
$ws.Cells.Item(493, 8).Value = "Sin especificar"
$ws.Cells.Item(493, 9).Value = "Primera"
$ws.Cells.Item(493, 10).Value = 400
$ws.Cells.Item(493, 11).Value = 10000
$ws.Cells.Item(493, 12).Value = 10000
$ws.Cells.Item(493, 13).Value = 10000
$ws.Cells.Item(493, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(493, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(493, 16).Value = 200
$ws.Cells.Item(493, 17).Value = 50
$ws.Cells.Item(493, 18).Value = "Hortaliza"
